$wb = $excel.ActiveWorkbook

$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet updates ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Foxleigh Coal Mine, Australia, M0040, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$newVersionText = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

foreach ($row in 2..8) {
    $wsData.Cells.Item($row, 19).Value = $newVersionText # column S = 19
}
